$wb = $excel.ActiveWorkbook

# ALC row 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 22828.5
$ws.Range("J3").Value = 22828.5
$ws.Range("L3").Value = 22828.5
$ws.Range("N3").Value = -23056.5

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 212
$ws.Range("I10").Value = 225
$ws.Range("J10").Value = 199
$ws.Range("K10").Value = 225
$ws.Range("L10").Value = 199
$ws.Range("M10").Value = 68
$ws.Range("N10").Value = -785

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 402.25
$ws.Range("I28").Value = 365.16666
$ws.Range("K28").Value = 365.16666
$ws.Range("M28").Value = 119.83334

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 259.35715
$ws.Range("I38").Value = 135.91667
$ws.Range("K38").Value = 407.75001
$ws.Range("M38").Value = -35.75001000000003

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 26.142857
$ws.Range("I39").Value = 22.666666
$ws.Range("J39").Value = 47
$ws.Range("K39").Value = 67.99999800000001
$ws.Range("L39").Value = 141
$ws.Range("M39").Value = 228.000002
$ws.Range("N39").Value = -733

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 475.25
$ws.Range("I41").Value = 567
$ws.Range("K41").Value = 567
$ws.Range("M41").Value = -127

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1963.7142
$ws.Range("I96").Value = 1963.7142
$ws.Range("K96").Value = 5891.142599999999
$ws.Range("M96").Value = -4518.142599999999

# ALC row 102
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 22828.5
$ws.Range("J102").Value = 22828.5
$ws.Range("L102").Value = 22828.5
$ws.Range("N102").Value = -29318.5

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8000
$ws.Range("I116").Value = 9000
$ws.Range("K116").Value = 9000
$ws.Range("M116").Value = -5558

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3666.76
$ws.Range("I138").Value = 4519
$ws.Range("J138").Value = 3098.6
$ws.Range("K138").Value = 13557
$ws.Range("L138").Value = 9295.799999999999
$ws.Range("M138").Value = -8417
$ws.Range("N138").Value = -19575.8

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3171.6667
$ws.Range("J74").Value = 3188.6667
$ws.Range("L74").Value = 3188.6667
$ws.Range("N74").Value = -4936.6667

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3171.6667
$ws.Range("J77").Value = 3188.6667
$ws.Range("L77").Value = 15943.3335
$ws.Range("N77").Value = -24679.3335

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5556263
$ws.Range("I7").Value = 10000422
$ws.Range("J7").Value = 1065
$ws.Range("K7").Value = 10000422
$ws.Range("L7").Value = 1065
$ws.Range("M7").Value = -10000309
$ws.Range("N7").Value = -1291

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1123.5
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 998
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 998
$ws.Range("M64").Value = -1275
$ws.Range("N64").Value = -1448

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1123.5
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 998
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 998
$ws.Range("M67").Value = -720
$ws.Range("N67").Value = -2558

# CRP row 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1930.3334
$ws.Range("I35").Value = 2417.5715
$ws.Range("K35").Value = 2417.5715
$ws.Range("M35").Value = -2123.5715

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 83495.8
$ws.Range("I62").Value = 2495
$ws.Range("J62").Value = 204997
$ws.Range("K62").Value = 2495
$ws.Range("L62").Value = 204997
$ws.Range("M62").Value = -1871
$ws.Range("N62").Value = -206245

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 83495.8
$ws.Range("I65").Value = 2495
$ws.Range("J65").Value = 204997
$ws.Range("K65").Value = 12475
$ws.Range("L65").Value = 1024985
$ws.Range("M65").Value = -9355
$ws.Range("N65").Value = -1031225

# CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61498

# CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 19499.857
$ws.Range("I69").Value = 17749.834
$ws.Range("K69").Value = 17749.834
$ws.Range("M69").Value = -17000.834

# CRP row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -187488

# CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 19499.857
$ws.Range("I72").Value = 17749.834
$ws.Range("K72").Value = 53249.50199999999
$ws.Range("M72").Value = -49505.50199999999

# CRP row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 36599.75
$ws.Range("J74").Value = 36599.75
$ws.Range("L74").Value = 36599.75
$ws.Range("N74").Value = -38347.75

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 36599.75
$ws.Range("J77").Value = 36599.75
$ws.Range("L77").Value = 109799.25
$ws.Range("N77").Value = -118535.25

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 40.92857
$ws.Range("I2").Value = 25.857143
$ws.Range("J2").Value = 56
$ws.Range("K2").Value = 155.142858
$ws.Range("L2").Value = 336
$ws.Range("M2").Value = -42.14285799999999
$ws.Range("N2").Value = -562

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 452.8
$ws.Range("I17").Value = 388.16666
$ws.Range("J17").Value = 549.75
$ws.Range("K17").Value = 1164.49998
$ws.Range("L17").Value = 1649.25
$ws.Range("M17").Value = -995.4999800000001
$ws.Range("N17").Value = -1987.25

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4916.6665
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 400
$ws.Range("I87").Value = 400
$ws.Range("K87").Value = 1200
$ws.Range("M87").Value = 48

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 400
$ws.Range("I90").Value = 400
$ws.Range("K90").Value = 3600
$ws.Range("M90").Value = 2640

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1873.4445
$ws.Range("I97").Value = 229.5
$ws.Range("J97").Value = 3188.6
$ws.Range("K97").Value = 688.5
$ws.Range("L97").Value = 9565.799999999999
$ws.Range("M97").Value = -192.5
$ws.Range("N97").Value = -10557.8

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = 2070

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 431.27777
$ws.Range("I2").Value = 146.2
$ws.Range("J2").Value = 787.625
$ws.Range("K2").Value = 146.2
$ws.Range("L2").Value = 787.625
$ws.Range("M2").Value = -33.19999999999999
$ws.Range("N2").Value = -1013.625

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1314.7142
$ws.Range("I14").Value = 240.6
$ws.Range("K14").Value = 240.6
$ws.Range("M14").Value = -72.59999999999999

# GSM row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 53199.6
$ws.Range("J63").Value = 53199.6
$ws.Range("L63").Value = 53199.6
$ws.Range("N63").Value = -54571.6

# GSM row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 53199.6
$ws.Range("J66").Value = 53199.6
$ws.Range("L66").Value = 159598.8
$ws.Range("N66").Value = -166462.8

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5945
$ws.Range("I7").Value = 5945
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5945
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5833
$ws.Range("N7").ClearContents()

# LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 26124.5
$ws.Range("I34").Value = 31666
$ws.Range("J34").Value = 9500
$ws.Range("K34").Value = 31666
$ws.Range("L34").Value = 9500
$ws.Range("M34").Value = -31494
$ws.Range("N34").Value = -9844

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3913.5
$ws.Range("I46").Value = 1383.3334
$ws.Range("J46").Value = 4997.857
$ws.Range("K46").Value = 1383.3334
$ws.Range("L46").Value = 4997.857
$ws.Range("M46").Value = -1195.3334
$ws.Range("N46").Value = -5373.857

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5945
$ws.Range("I126").Value = 5945
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 17835
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -15365
$ws.Range("N126").ClearContents()

# WVR row 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 42464
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 42464
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 42464
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -42762

# WVR row 103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 13463
$ws.Range("J103").Value = 13463
$ws.Range("L103").Value = 13463
$ws.Range("N103").Value = -15807

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1189.2307
$ws.Range("I126").Value = 1213.3334
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 3640.0002
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -1170.0002
$ws.Range("N126").Value = -7640
